# [#571] Update conversion rate
# conversion_rates.xlsx — fill in placeholder/derived conversion rates for
# Zimbabwe on both the "USD conversion" and "EUR conversion" sheets, and
# restate the EUR "USD conversion" Ecuador precision artifact that came
# along with the copy.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("USD conversion")
$ws2 = $wb.Worksheets.Item("EUR conversion")

# ---------------------------------------------------------------------
# USD conversion — Zimbabwe (row 59) gets a flat placeholder rate of 1
# across every month column (B:O).
# ---------------------------------------------------------------------
$ws1.Range("B59:O59").Value = 1

# ---------------------------------------------------------------------
# EUR conversion — Ecuador (row 19) value restated with trimmed
# precision, and Zimbabwe (row 58) populated as a copy of Ecuador's
# conversion rates (B:O).
# ---------------------------------------------------------------------
$ws2.Range("B19").Value = 1.3257166666667

$ws2.Range("B58").Value = 1.3257166666667
$ws2.Range("C58").Value = 1.3919552529182999
$ws2.Range("D58").Value = 1.2847886718749999
$ws2.Range("E58").Value = 1.328118039215701
$ws2.Range("F58").Value = 1.3285007843137
$ws2.Range("G58").Value = 1.1095128906249989
$ws2.Range("H58").Value = 1.1069031128405
$ws2.Range("I58").Value = 1.1296811764706001
$ws2.Range("J58").Value = 1.1809545098039
$ws2.Range("K58").Value = 1.1194745098039001
$ws2.Range("L58").Value = 1.1421961089493999
$ws2.Range("M58").Value = 1.1827403100775
$ws2.Range("N58").Value = 1.0530486381323001
$ws2.Range("O58").Value = 1.0812686274509991

# ---------------------------------------------------------------------
# View / selection state.
# Before: "USD conversion" tab selected, selection R55 on sheet1,
#         selection O50 on sheet2 (scrolled to A31).
# After:  "EUR conversion" tab selected (now active), scrolled to A28,
#         selection B58:O58 on sheet2, selection L66 on sheet1.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("L66").Select()

$ws2.Activate()
$ws2.Range("A28").Select()
$ws2.Range("B58:O58").Select()
